$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Locate the row whose first cell contains "fabric" (exact match, ignoring the
# trailing cell-mark characters Word appends to cell Range.Text).
$fabricRow = -1
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    $txt = $t.Cell($i, 1).Range.Text
    $txt = $txt.TrimEnd([char]13, [char]7)
    if ($txt -eq "fabric") {
        $fabricRow = $i
        break
    }
}

if ($fabricRow -gt 0) {
    # Insert a brand-new row immediately above the "fabric" row and fill it in
    # with the "demjson" package entry.
    $newRow = $t.Rows.Add($t.Rows.Item($fabricRow))

    # The newly inserted row copies paragraph formatting (including any
    # explicit tab stops) from the neighboring "cryptography" row; clear that
    # out so the new row's paragraphs stay plain, matching a freshly typed row.
    $t.Cell($fabricRow, 1).Range.ParagraphFormat.TabStops.ClearAll()
    $t.Cell($fabricRow, 2).Range.ParagraphFormat.TabStops.ClearAll()

    $t.Cell($fabricRow, 1).Range.Text = "demjson"
    $t.Cell($fabricRow, 2).Range.Text = "A library for JSON format data check"
}
